# Refresh the "cryptos" price list with the latest scraped values.
# (GitHub Actions data-refresh commit: "Updated cryptos list ... with GitHub Actions")
#
# Price (col D) and Volume(1h) (col E) cells are plain text in this sheet
# (e.g. "58.487.91", "  +0.29%  "), so values that otherwise look like plain
# numbers/decimals are written with a leading apostrophe to force Excel to
# keep storing them as text instead of silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range('D2').Value = '58.487.91'

$ws.Range('D3').Value = '3.097.47'
$ws.Range('E3').Value = '  +0.29%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''528.09'
$ws.Range('E5').Value = '  +2.29%  '

$ws.Range('D6').Value = '''142.79'
$ws.Range('E6').Value = '  +0.94%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').Value = '  +1.51%  '

$ws.Range('E9').Value = '  +0.78%  '

$ws.Range('E10').Value = '  +0.67%  '

$ws.Range('E11').Value = '  +2.63%  '

$ws.Range('D12').Value = '3.626.29'
$ws.Range('E12').Value = '  +0.33%  '

$ws.Range('E13').Value = '  +0.87%  '

$ws.Range('D14').Value = '''26.86'
$ws.Range('E14').Value = '  +4.83%  '

$ws.Range('E15').Value = '  +1.84%  '

$ws.Range('D16').Value = '58.518.34'
$ws.Range('E16').Value = '  +1.67%  '

$ws.Range('D17').Value = '3.108.99'
$ws.Range('E17').Value = '  +0.64%  '

$ws.Range('E18').Value = '  -0.51%  '

$ws.Range('D19').Value = '''12.91'
$ws.Range('E19').Value = '  -2.21%  '

$ws.Range('D20').Value = '''8.08'
$ws.Range('E20').Value = '  -0.92%  '

$ws.Range('D21').Value = '''341.78'
$ws.Range('E21').Value = '  +2.13%  '

$ws.Range('D23').Value = '''0.505'
$ws.Range('E23').Value = '  +0.80%  '

$ws.Range('D24').Value = '''66.00'
$ws.Range('E24').Value = '  +0.15%  '

$ws.Range('E25').Value = '  +0.22%  '

$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('E27').Value = '  +0.30%  '

$ws.Range('D28').Value = '''6.63'
$ws.Range('E28').Value = '  +3.53%  '

$ws.Range('E29').Value = '  +1.50%  '

$ws.Range('E30').Value = '  +3.06%  '

$ws.Range('D31').Value = '''1.21'
$ws.Range('E31').Value = '  +3.59%  '

$ws.Range('D32').Value = '''20.97'
$ws.Range('E32').Value = '  +0.55%  '

$ws.Range('D33').Value = '''154.09'
$ws.Range('E33').Value = '  -0.51%  '

$ws.Range('D34').Value = '''4.66'
$ws.Range('E34').Value = '  +2.58%  '

$ws.Range('D35').Value = '''6.07'
$ws.Range('E35').Value = '  +3.07%  '

$ws.Range('D36').Value = '''26.88'
$ws.Range('E36').Value = '  -3.86%  '

$ws.Range('E37').Value = '  +3.20%  '

$ws.Range('E38').Value = '  +0.31%  '

$ws.Range('D39').Value = '3.137.19'
$ws.Range('E39').Value = '  +0.30%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '''3.88'
$ws.Range('E40').Value = '  +0.28%  '

$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '''0.677'
$ws.Range('E41').Value = '  +0.89%  '

$ws.Range('D42').Value = '''36.78'
$ws.Range('E42').Value = '  -0.03%  '

$ws.Range('E43').Value = '  +7.90%  '

$ws.Range('D44').Value = '''0.999'
$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('D45').Value = '2.296.27'
$ws.Range('E45').Value = '  +0.27%  '

$ws.Range('E46').Value = '  +0.63%  '

$ws.Range('D47').Value = '''20.84'
$ws.Range('E47').Value = '  +4.15%  '

$ws.Range('D48').Value = '''0.962'
$ws.Range('E48').Value = '  +2.27%  '

$ws.Range('D49').Value = '''5.99'
$ws.Range('E49').Value = '  +1.72%  '

$ws.Range('D50').Value = '''267.58'
$ws.Range('E50').Value = '  +5.86%  '

$ws.Range('D51').Value = '''0.743'
$ws.Range('E51').Value = '  +8.12%  '
